# Weekly update for Hortaliza, Femacal de La Calera - Arveja Verde.
# A new daily price record (44677) is inserted as row 10, pushing every
# existing record (previously rows 10-59) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 10 (shifts rows 10..59 down to 11..60).
$ws.Rows(10).Insert()

# Populate the newly inserted row with the new market record.
$ws.Cells.Item(10, 1).Value  = 3
$ws.Cells.Item(10, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(10, 3).Value  = "Coquimbo"
$ws.Cells.Item(10, 4).Value  = 44677
$ws.Cells.Item(10, 5).Value  = 5
$ws.Cells.Item(10, 6).Value  = 100112022
$ws.Cells.Item(10, 7).Value  = "Arveja Verde"
$ws.Cells.Item(10, 8).Value  = "Perfection"
$ws.Cells.Item(10, 9).Value  = "Primera"
$ws.Cells.Item(10, 10).Value = 65
$ws.Cells.Item(10, 11).Value = 22000
$ws.Cells.Item(10, 12).Value = 23000
$ws.Cells.Item(10, 13).Value = 22462
$ws.Cells.Item(10, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(10, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 16).Value = 898
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"
